$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"27.13486633333333"
$ws.Range("H2").Value = [double]"81.40459899999999"
$ws.Range("I2").Value = [double]"0.04747038381101173"
$ws.Range("J2").Value = [double]"0.05206311700485852"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"0.2524303333333334"
$ws.Range("N2").Value = [double]"0.757291"
$ws.Range("O2").Value = [double]"0.0004152138066912919"
$ws.Range("P2").Value = [double]"0.0004153435393530768"
$ws.Range("Q2").Value = [double]"6.849663353478778"
$ws.Range("R2").Value = [double]"61.64697018130899"
$ws.Range("S2").Value = [double]"1.971035876726686E-05"
$ws.Range("T2").Value = [double]"2.16240792865513E-05"

$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"27.13486633333333"
$ws.Range("H3").Value = [double]"81.40459899999999"
$ws.Range("I3").Value = [double]"0.04747038381101173"
$ws.Range("J3").Value = [double]"0.05206311700485852"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"0.3565143333333333"
$ws.Range("N3").Value = [double]"1.069543"
$ws.Range("O3").Value = [double]"0.0005864179297654724"
$ws.Range("P3").Value = [double]"0.0005866011547876679"
$ws.Range("Q3").Value = [double]"9.673968780917443"
$ws.Range("R3").Value = [double]"87.06571902825698"
$ws.Range("S3").Value = [double]"2.783748419962589E-05"
$ws.Range("T3").Value = [double]"3.054028455689547E-05"

$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"27.13486633333333"
$ws.Range("H4").Value = [double]"81.40459899999999"
$ws.Range("I4").Value = [double]"0.04747038381101173"
$ws.Range("J4").Value = [double]"0.05206311700485852"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"343.6225723333334"
$ws.Range("N4").Value = [double]"1030.867717"
$ws.Range("O4").Value = [double]"0.5652127240000627"
$ws.Range("P4").Value = [double]"0.5653893235012776"
$ws.Range("Q4").Value = [double]"9324.152569381165"
$ws.Range("R4").Value = [double]"83917.37312443048"
$ws.Range("S4").Value = [double]"0.02683086494315042"
$ws.Range("T4").Value = [double]"0.02943593050274482"

$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"27.13486633333333"
$ws.Range("H5").Value = [double]"81.40459899999999"
$ws.Range("I5").Value = [double]"0.04747038381101173"
$ws.Range("J5").Value = [double]"0.05206311700485852"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"263.1514486666667"
$ws.Range("N5").Value = [double]"789.454346"
$ws.Range("O5").Value = [double]"0.432848593488691"
$ws.Range("P5").Value = [double]"0.4329838361017211"
$ws.Range("Q5").Value = [double]"7140.579384993028"
$ws.Range("R5").Value = [double]"64265.21446493724"
$ws.Range("S5").Value = [double]"0.02054748886496475"
$ws.Range("T5").Value = [double]"0.02254248812017639"

$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"27.13486633333333"
$ws.Range("H6").Value = [double]"81.40459899999999"
$ws.Range("I6").Value = [double]"0.04747038381101173"
$ws.Range("J6").Value = [double]"0.05206311700485852"
$ws.Range("K6").Value = [double]"2"
$ws.Range("M6").Value = [double]"0.5696825"
$ws.Range("N6").Value = [double]"1.139365"
$ws.Range("O6").Value = [double]"0.000937050774789659"
$ws.Range("P6").Value = [double]"0.0006248957028606154"
$ws.Range("Q6").Value = [double]"15.45825848993917"
$ws.Range("R6").Value = [double]"92.74955093963499"
$ws.Range("S6").Value = [double]"4.448215992967103E-05"
$ws.Range("T6").Value = [double]"3.253401809386552E-05"

$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"207.121208"
$ws.Range("H7").Value = [double]"621.3636240000001"
$ws.Range("I7").Value = [double]"0.3623427924198875"
$ws.Range("J7").Value = [double]"0.3973992557702412"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"0.2524303333333334"
$ws.Range("N7").Value = [double]"0.757291"
$ws.Range("O7").Value = [double]"0.0004152138066912919"
$ws.Range("P7").Value = [double]"0.0004153435393530768"
$ws.Range("Q7").Value = [double]"52.28367557584268"
$ws.Range("R7").Value = [double]"470.5530801825841"
$ws.Range("S7").Value = [double]"0.0001504497301678141"
$ws.Range("T7").Value = [double]"0.0001650572134278906"

$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"207.121208"
$ws.Range("H8").Value = [double]"621.3636240000001"
$ws.Range("I8").Value = [double]"0.3623427924198875"
$ws.Range("J8").Value = [double]"0.3973992557702412"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"0.3565143333333333"
$ws.Range("N8").Value = [double]"1.069543"
$ws.Range("O8").Value = [double]"0.0005864179297654724"
$ws.Range("P8").Value = [double]"0.0005866011547876679"
$ws.Range("Q8").Value = [double]"73.84167938931468"
$ws.Range("R8").Value = [double]"664.575114503832"
$ws.Range("S8").Value = [double]"0.0002124843101963108"
$ws.Range("T8").Value = [double]"0.0002331148623465833"

$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"207.121208"
$ws.Range("H9").Value = [double]"621.3636240000001"
$ws.Range("I9").Value = [double]"0.3623427924198875"
$ws.Range("J9").Value = [double]"0.3973992557702412"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"343.6225723333334"
$ws.Range("N9").Value = [double]"1030.867717"
$ws.Range("O9").Value = [double]"0.5652127240000627"
$ws.Range("P9").Value = [double]"0.5653893235012776"
$ws.Range("Q9").Value = [double]"71171.5222777474"
$ws.Range("R9").Value = [double]"640543.7004997266"
$ws.Range("S9").Value = [double]"0.2048007567254339"
$ws.Range("T9").Value = [double]"0.2246852963798479"

$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"207.121208"
$ws.Range("H10").Value = [double]"621.3636240000001"
$ws.Range("I10").Value = [double]"0.3623427924198875"
$ws.Range("J10").Value = [double]"0.3973992557702412"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"263.1514486666667"
$ws.Range("N10").Value = [double]"789.454346"
$ws.Range("O10").Value = [double]"0.432848593488691"
$ws.Range("P10").Value = [double]"0.4329838361017211"
$ws.Range("Q10").Value = [double]"54504.24593479"
$ws.Range("R10").Value = [double]"490538.2134131099"
$ws.Range("S10").Value = [double]"0.1568395680597131"
$ws.Range("T10").Value = [double]"0.1720674542273681"

$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"207.121208"
$ws.Range("H11").Value = [double]"621.3636240000001"
$ws.Range("I11").Value = [double]"0.3623427924198875"
$ws.Range("J11").Value = [double]"0.3973992557702412"
$ws.Range("K11").Value = [double]"2"
$ws.Range("M11").Value = [double]"0.5696825"
$ws.Range("N11").Value = [double]"1.139365"
$ws.Range("O11").Value = [double]"0.000937050774789659"
$ws.Range("P11").Value = [double]"0.0006248957028606154"
$ws.Range("Q11").Value = [double]"117.99332757646"
$ws.Range("R11").Value = [double]"707.9599654587601"
$ws.Range("S11").Value = [double]"0.0003395335943765042"
$ws.Range("T11").Value = [double]"0.0002483330872508304"

$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"84.750407"
$ws.Range("H12").Value = [double]"254.251221"
$ws.Range("I12").Value = [double]"0.1482643879283573"
$ws.Range("J12").Value = [double]"0.1626088848807073"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"0.2524303333333334"
$ws.Range("N12").Value = [double]"0.757291"
$ws.Range("O12").Value = [double]"0.0004152138066912919"
$ws.Range("P12").Value = [double]"0.0004153435393530768"
$ws.Range("Q12").Value = [double]"21.39357348914567"
$ws.Range("R12").Value = [double]"192.542161402311"
$ws.Range("S12").Value = [double]"6.156142090848766E-05"
$ws.Range("T12").Value = [double]"6.753854977661E-05"

$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"84.750407"
$ws.Range("H13").Value = [double]"254.251221"
$ws.Range("I13").Value = [double]"0.1482643879283573"
$ws.Range("J13").Value = [double]"0.1626088848807073"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"0.3565143333333333"
$ws.Range("N13").Value = [double]"1.069543"
$ws.Range("O13").Value = [double]"0.0005864179297654724"
$ws.Range("P13").Value = [double]"0.0005866011547876679"
$ws.Range("Q13").Value = [double]"30.21473485133366"
$ws.Range("R13").Value = [double]"271.932613662003"
$ws.Range("S13").Value = [double]"8.694489542689217E-05"
$ws.Range("T13").Value = [double]"9.538655964975786E-05"

$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"84.750407"
$ws.Range("H14").Value = [double]"254.251221"
$ws.Range("I14").Value = [double]"0.1482643879283573"
$ws.Range("J14").Value = [double]"0.1626088848807073"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"343.6225723333334"
$ws.Range("N14").Value = [double]"1030.867717"
$ws.Range("O14").Value = [double]"0.5652127240000627"
$ws.Range("P14").Value = [double]"0.5653893235012776"
$ws.Range("Q14").Value = [double]"29122.15285963694"
$ws.Range("R14").Value = [double]"262099.3757367325"
$ws.Range("S14").Value = [double]"0.08380091857318883"
$ws.Range("T14").Value = [double]"0.09193732741800023"

$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"84.750407"
$ws.Range("H15").Value = [double]"254.251221"
$ws.Range("I15").Value = [double]"0.1482643879283573"
$ws.Range("J15").Value = [double]"0.1626088848807073"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"263.1514486666667"
$ws.Range("N15").Value = [double]"789.454346"
$ws.Range("O15").Value = [double]"0.432848593488691"
$ws.Range("P15").Value = [double]"0.4329838361017211"
$ws.Range("Q15").Value = [double]"22302.19237713961"
$ws.Range("R15").Value = [double]"200719.7313942565"
$ws.Range("S15").Value = [double]"0.0641760317792511"
$ws.Range("T15").Value = [double]"0.07040701875987182"

$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"84.750407"
$ws.Range("H16").Value = [double]"254.251221"
$ws.Range("I16").Value = [double]"0.1482643879283573"
$ws.Range("J16").Value = [double]"0.1626088848807073"
$ws.Range("K16").Value = [double]"2"
$ws.Range("M16").Value = [double]"0.5696825"
$ws.Range("N16").Value = [double]"1.139365"
$ws.Range("O16").Value = [double]"0.000937050774789659"
$ws.Range("P16").Value = [double]"0.0006248957028606154"
$ws.Range("Q16").Value = [double]"48.2808237357775"
$ws.Range("R16").Value = [double]"289.684942414665"
$ws.Range("S16").Value = [double]"0.0001389312595819818"
$ws.Range("T16").Value = [double]"0.0001016135934089105"

$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"101.3352343333333"
$ws.Range("H17").Value = [double]"304.005703"
$ws.Range("I17").Value = [double]"0.1772782813185584"
$ws.Range("J17").Value = [double]"0.1944298562963656"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"0.2524303333333334"
$ws.Range("N17").Value = [double]"0.757291"
$ws.Range("O17").Value = [double]"0.0004152138066912919"
$ws.Range("P17").Value = [double]"0.0004153435393530768"
$ws.Range("Q17").Value = [double]"25.58008698117478"
$ws.Range("R17").Value = [double]"230.220782830573"
$ws.Range("S17").Value = [double]"7.360839002996838E-05"
$ws.Range("T17").Value = [double]"8.07551846700426E-05"

$ws.Range("E18").Value = [double]"3"
$ws.Range("G18").Value = [double]"101.3352343333333"
$ws.Range("H18").Value = [double]"304.005703"
$ws.Range("I18").Value = [double]"0.1772782813185584"
$ws.Range("J18").Value = [double]"0.1944298562963656"
$ws.Range("K18").Value = [double]"3"
$ws.Range("M18").Value = [double]"0.3565143333333333"
$ws.Range("N18").Value = [double]"1.069543"
$ws.Range("O18").Value = [double]"0.0005864179297654724"
$ws.Range("P18").Value = [double]"0.0005866011547876679"
$ws.Range("Q18").Value = [double]"36.12746351152544"
$ws.Range("R18").Value = [double]"325.147171603729"
$ws.Range("S18").Value = [double]"0.00010395916272321"
$ws.Range("T18").Value = [double]"0.0001140527782286484"

$ws.Range("E19").Value = [double]"3"
$ws.Range("G19").Value = [double]"101.3352343333333"
$ws.Range("H19").Value = [double]"304.005703"
$ws.Range("I19").Value = [double]"0.1772782813185584"
$ws.Range("J19").Value = [double]"0.1944298562963656"
$ws.Range("K19").Value = [double]"3"
$ws.Range("M19").Value = [double]"343.6225723333334"
$ws.Range("N19").Value = [double]"1030.867717"
$ws.Range("O19").Value = [double]"0.5652127240000627"
$ws.Range("P19").Value = [double]"0.5653893235012776"
$ws.Range("Q19").Value = [double]"34821.07388962112"
$ws.Range("R19").Value = [double]"313389.6650065901"
$ws.Range("S19").Value = [double]"0.1001999402901118"
$ws.Range("T19").Value = [double]"0.1099285649198528"

$ws.Range("E20").Value = [double]"3"
$ws.Range("G20").Value = [double]"101.3352343333333"
$ws.Range("H20").Value = [double]"304.005703"
$ws.Range("I20").Value = [double]"0.1772782813185584"
$ws.Range("J20").Value = [double]"0.1944298562963656"
$ws.Range("K20").Value = [double]"3"
$ws.Range("M20").Value = [double]"263.1514486666667"
$ws.Range("N20").Value = [double]"789.454346"
$ws.Range("O20").Value = [double]"0.432848593488691"
$ws.Range("P20").Value = [double]"0.4329838361017211"
$ws.Range("Q20").Value = [double]"26666.5137157928"
$ws.Range("R20").Value = [double]"239998.6234421352"
$ws.Range("S20").Value = [double]"0.07673465472483049"
$ws.Range("T20").Value = [double]"0.08418498503190676"

$ws.Range("E21").Value = [double]"3"
$ws.Range("G21").Value = [double]"101.3352343333333"
$ws.Range("H21").Value = [double]"304.005703"
$ws.Range("I21").Value = [double]"0.1772782813185584"
$ws.Range("J21").Value = [double]"0.1944298562963656"
$ws.Range("K21").Value = [double]"2"
$ws.Range("M21").Value = [double]"0.5696825"
$ws.Range("N21").Value = [double]"1.139365"
$ws.Range("O21").Value = [double]"0.000937050774789659"
$ws.Range("P21").Value = [double]"0.0006248957028606154"
$ws.Range("Q21").Value = [double]"57.72890963309916"
$ws.Range("R21").Value = [double]"346.373457798595"
$ws.Range("S21").Value = [double]"0.0001661187508629343"
$ws.Range("T21").Value = [double]"0.0001214983817074058"

$ws.Range("E22").Value = [double]"2"
$ws.Range("G22").Value = [double]"151.2750305"
$ws.Range("H22").Value = [double]"302.550061"
$ws.Range("I22").Value = [double]"0.2646441545221851"
$ws.Range("J22").Value = [double]"0.1934988860478274"
$ws.Range("K22").Value = [double]"3"
$ws.Range("M22").Value = [double]"0.2524303333333334"
$ws.Range("N22").Value = [double]"0.757291"
$ws.Range("O22").Value = [double]"0.0004152138066912919"
$ws.Range("P22").Value = [double]"0.0004153435393530768"
$ws.Range("Q22").Value = [double]"38.18640637412517"
$ws.Range("R22").Value = [double]"229.118438244751"
$ws.Range("S22").Value = [double]"0.000109883906817755"
$ws.Range("T22").Value = [double]"8.036851219198231E-05"

$ws.Range("E23").Value = [double]"2"
$ws.Range("G23").Value = [double]"151.2750305"
$ws.Range("H23").Value = [double]"302.550061"
$ws.Range("I23").Value = [double]"0.2646441545221851"
$ws.Range("J23").Value = [double]"0.1934988860478274"
$ws.Range("K23").Value = [double]"3"
$ws.Range("M23").Value = [double]"0.3565143333333333"
$ws.Range("N23").Value = [double]"1.069543"
$ws.Range("O23").Value = [double]"0.0005864179297654724"
$ws.Range("P23").Value = [double]"0.0005866011547876679"
$ws.Range("Q23").Value = [double]"53.93171664868716"
$ws.Range("R23").Value = [double]"323.5902998921229"
$ws.Range("S23").Value = [double]"0.0001551920772194336"
$ws.Range("T23").Value = [double]"0.0001135066700057829"

$ws.Range("E24").Value = [double]"2"
$ws.Range("G24").Value = [double]"151.2750305"
$ws.Range("H24").Value = [double]"302.550061"
$ws.Range("I24").Value = [double]"0.2646441545221851"
$ws.Range("J24").Value = [double]"0.1934988860478274"
$ws.Range("K24").Value = [double]"3"
$ws.Range("M24").Value = [double]"343.6225723333334"
$ws.Range("N24").Value = [double]"1030.867717"
$ws.Range("O24").Value = [double]"0.5652127240000627"
$ws.Range("P24").Value = [double]"0.5653893235012776"
$ws.Range("Q24").Value = [double]"51981.51511021346"
$ws.Range("R24").Value = [double]"311889.0906612807"
$ws.Range("S24").Value = [double]"0.1495802434681777"
$ws.Range("T24").Value = [double]"0.1094022042808319"

$ws.Range("E25").Value = [double]"2"
$ws.Range("G25").Value = [double]"151.2750305"
$ws.Range("H25").Value = [double]"302.550061"
$ws.Range("I25").Value = [double]"0.2646441545221851"
$ws.Range("J25").Value = [double]"0.1934988860478274"
$ws.Range("K25").Value = [double]"3"
$ws.Range("M25").Value = [double]"263.1514486666667"
$ws.Range("N25").Value = [double]"789.454346"
$ws.Range("O25").Value = [double]"0.432848593488691"
$ws.Range("P25").Value = [double]"0.4329838361017211"
$ws.Range("Q25").Value = [double]"39808.24342316918"
$ws.Range("R25").Value = [double]"238849.4605390151"
$ws.Range("S25").Value = [double]"0.1145508500599316"
$ws.Range("T25").Value = [double]"0.0837818899623981"

$ws.Range("E26").Value = [double]"2"
$ws.Range("G26").Value = [double]"151.2750305"
$ws.Range("H26").Value = [double]"302.550061"
$ws.Range("I26").Value = [double]"0.2646441545221851"
$ws.Range("J26").Value = [double]"0.1934988860478274"
$ws.Range("K26").Value = [double]"2"
$ws.Range("M26").Value = [double]"0.5696825"
$ws.Range("N26").Value = [double]"1.139365"
$ws.Range("O26").Value = [double]"0.000937050774789659"
$ws.Range("P26").Value = [double]"0.0006248957028606154"
$ws.Range("Q26").Value = [double]"86.17873756281624"
$ws.Range("R26").Value = [double]"344.714950251265"
$ws.Range("S26").Value = [double]"0.0002479850100385678"
$ws.Range("T26").Value = [double]"0.0001209166223996032"

